# Expedia_automation/Data_sheets/Data_sheets_expedia.xlsx
# Commit: "Added some instructions to the excel sheet."
#
# Summary of the edit:
#  - FlightsOnlyPositive (sheet1): lowercase a couple of the sample city
#    values, fix a stale departure date, widen two spare columns and drop
#    an instructions note (merged, small wrapped font) into G12:H13.
#  - FlightsWithHotelsPositive (sheet2): becomes the active tab, one city
#    value is tidied up, the hotel check-out date is re-typed as text, a
#    brand-new sample row (New York / Delhi) is added, and the same kind
#    of instructions note is dropped into I11:J12.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$instructions = "Instruction - Put the start tag (One-way) before the first column and the first row to be used and the end tag (One-way) after the last column and last row to be used.`nOnly use future dates here."

# ---------------------------------------------------------------------
# 1. FlightsWithHotelsPositive (sheet2) data edits - done first so the
#    shared-string table fills up in the same order the original commit
#    produced.
# ---------------------------------------------------------------------

# Hotel check-out date on the existing sample row was a raw date serial;
# re-enter it as plain text.
$ws2.Range("F3").NumberFormat = "@"
$ws2.Range("F3").Font.Name = "Calibri"
$ws2.Range("F3").Value = "18/12/2019"

# Insert a brand-new sample row (pushes the trailing "OneWayHotels" end
# tag from G4 down to G5, matching the target layout).
$ws2.Rows.Item(4).Insert()
$ws2.Range("B4").Value = "new york"

# ---------------------------------------------------------------------
# 2. FlightsOnlyPositive (sheet1) data edits.
# ---------------------------------------------------------------------
$ws1.Range("C3").Value = "hyderabad"
$ws1.Range("C6").Value = "mumbai"

# ---------------------------------------------------------------------
# back to sheet2 for the rest of the new row + the origin city tidy-up
# ---------------------------------------------------------------------
$ws2.Range("C4").Value = "delhi"

$ws1.Range("B4").Value = "delhi"
$ws1.Range("B6").Value = "Delhi"

$ws2.Range("B2").Value = "Delhi"

$ws1.Range("D6").Value = "09/01/2020"

$ws2.Range("D4").Value = "07/12/2019"
$ws2.Range("E4").Value = "16/12/2019"
$ws2.Range("F4").NumberFormat = "@"
$ws2.Range("F4").Font.Name = "Calibri"
$ws2.Range("F4").Value = "18/12/2019"

# ---------------------------------------------------------------------
# 3. Formatting touch-ups on sheet1.
# ---------------------------------------------------------------------
$ws1.Rows.Item(2).RowHeight = 17

$ws1.Columns.Item(7).ColumnWidth = 23.714285714285715
$ws1.Columns.Item(8).ColumnWidth = 25.571428571428573

$ws1.Rows.Item(12).RowHeight = 31.5
$ws1.Rows.Item(13).RowHeight = 46

$rng1 = $ws1.Range("G12:H13")
$rng1.Merge()
$rng1.Font.Name = "Arial"
$rng1.Font.Size = 8
$rng1.VerticalAlignment = -4160
$rng1.WrapText = $true
$ws1.Range("G12").Value = $instructions

# ---------------------------------------------------------------------
# 4. Formatting touch-ups on sheet2.
# ---------------------------------------------------------------------
$ws2.Columns.Item(9).ColumnWidth = 19.142857142857142
$ws2.Columns.Item(10).ColumnWidth = 18.142857142857142

$ws2.Rows.Item(11).RowHeight = 49.5
$ws2.Rows.Item(12).RowHeight = 34

$rng2 = $ws2.Range("I11:J12")
$rng2.Merge()
$rng2.Font.Name = "Arial"
$rng2.Font.Size = 8
$rng2.VerticalAlignment = -4160
$rng2.WrapText = $true
$ws2.Range("I11").Value = $instructions

# ---------------------------------------------------------------------
# 5. Make FlightsWithHotelsPositive the active/visible tab, like in the
#    target workbook (activeTab="1" / tabSelected moves to sheet2).
# ---------------------------------------------------------------------
$ws2.Activate()

Write-Output "edit applied"
